# "add monte_carlo and update database"
# Fills in the quarterly income-statement figures (rows 11-27, columns D..M)
# that were previously placeholder zeros / "-" strings, per the new data
# pulled from the (fictitious) monte_carlo database update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => values for quarters D..M (10 columns)
$rows = @{
    11 = @(37415997, 63734615, 65227719, 78547663, 82700181, 100241751, 132230271, 136993252, 103239511, 82355074)
    12 = @(-36302111, -46820037, -55450974, -66365872, -76509899, -95791807, -109406194, -120997090, -100898806, -84576589)
    13 = @(1113886, 16914578, 9776745, 12181791, 6190282, 4449944, 22824077, 15996162, 2340705, -2221515)
    14 = @(-318950, -188291, -400738, -293345, -265388, -1233474, -924523, -287482, -561619, -1462162)
    15 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    16 = @(2939744, -1664696, 970463, -14819, -112774, 775536, 0, 1329369, 949847, 3449894)
    17 = @(3734680, 15061591, 10346470, 11873627, 5812120, 3992006, 21899554, 17038049, 2728933, -233783)
    18 = @(-51252, -249840, -297126, -220130, -318280, -247658, -266541, -250016, -362031, -413963)
    19 = @(-11726, -324413, 33240, 171331, -204571, 424271, 40068, 231596, 265952, 5768726)
    20 = @(3671702, 14487338, 10082584, 11824828, 5289269, 4168619, 21673081, 17019629, 2632854, 5120980)
    21 = @(-587934, -1479544, -2268582, -2645432, -1205239, 3286248, -5418270, -2892738, 1802232, 2028757)
    22 = @(3083768, 13007794, 7814002, 9179396, 4084030, 7454867, 16254811, 14126891, 4435086, 7149737)
    23 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    24 = @(3083768, 13007794, 7814002, 9179396, 4084030, 7454867, 16254811, 14126891, 4435086, 7149737)
    25 = @(2697, 1626, 977, 1147, 511, 6520, 2032, 1766, 554, 894)
    26 = @(1143422, 8000000, 8000000, 8000000, 8000000, 1143422, 8000000, 8000000, 8000000, 8000000)
    27 = @(385, 1626, 977, 1147, 511, 932, 2032, 1766, 554, 894)
}

# Columns D (4) through M (13)
$firstCol = 4

foreach ($r in $rows.Keys) {
    $values = $rows[$r]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($r, $firstCol + $i).Value = $values[$i]
    }
}
